# Apply updated crypto price/volume data to the worksheet.
# (Matches the GitHub Actions data-refresh commit on Mon Nov  6 23:56:36 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value to a cell, forcing text interpretation (via a leading
# apostrophe, same as typing it in the Excel UI) so numeric-looking text
# (e.g. trailing zeros like "75.00") is not silently converted to a number.
function Set-TextValue($address, $text) {
    $ws.Range($address).Value = '''' + $text
}

$ws.Range('D2').Value = '35.317.85'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '1.912.67'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  +0.04%  '
Set-TextValue 'D5' '0.725'
$ws.Range('E5').Value = '  +8.62%  '
Set-TextValue 'D6' '255.30'
$ws.Range('E6').Value = '  +3.61%  '
$ws.Range('E7').Value = '  +0.09%  '
Set-TextValue 'D8' '42.38'
$ws.Range('E8').Value = '  +1.87%  '
$ws.Range('E9').Value = '  +5.84%  '
$ws.Range('E10').Value = '  +0.38%  '
Set-TextValue 'D11' '0.0771'
$ws.Range('E11').Value = '  +7.29%  '
Set-TextValue 'D12' '0.0987'
$ws.Range('E12').Value = '  -0.48%  '
Set-TextValue 'D13' '13.10'
$ws.Range('E13').Value = '  +6.34%  '
$ws.Range('D14').Value = '2.189.55'
$ws.Range('E14').Value = '  +0.19%  '
Set-TextValue 'D15' '0.739'
$ws.Range('E15').Value = '  +5.70%  '
Set-TextValue 'D16' '5.01'
$ws.Range('E16').Value = '  +3.96%  '
$ws.Range('D17').Value = '1.904.49'
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('D18').Value = '35.317.75'
$ws.Range('E18').Value = '  -0.14%  '
Set-TextValue 'D19' '75.00'
$ws.Range('E19').Value = '  +3.93%  '
$ws.Range('D20').Value = '0.0₃0849'
$ws.Range('E20').Value = '  +3.51%  '
Set-TextValue 'D21' '245.61'
$ws.Range('E21').Value = '  +1.96%  '
Set-TextValue 'D22' '13.17'
$ws.Range('E22').Value = '  +5.57%  '
Set-TextValue 'D23' '5.16'
$ws.Range('E23').Value = '  +6.99%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('E25').Value = '  +7.32%  '
Set-TextValue 'D26' '2.39'
$ws.Range('E26').Value = '  -0.37%  '
Set-TextValue 'D27' '166.89'
$ws.Range('E27').Value = '  -2.40%  '
Set-TextValue 'D28' '8.80'
$ws.Range('E28').Value = '  +4.11%  '
Set-TextValue 'D29' '18.84'
$ws.Range('E29').Value = '  +2.69%  '
Set-TextValue 'D30' '0.132'
$ws.Range('E30').Value = '  +4.23%  '
$ws.Range('D31').Value = '4.128.81'
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('E32').Value = '  +27.03%  '
Set-TextValue 'D33' '4.36'
$ws.Range('E33').Value = '  +5.20%  '
$ws.Range('E34').Value = '  +15.32%  '
$ws.Range('E35').Value = '  +4.92%  '
$ws.Range('E36').Value = '  +4.55%  '
$ws.Range('E37').Value = '  +0.05%  '
Set-TextValue 'D38' '0.929'
$ws.Range('E38').Value = '  -1.86%  '
$ws.Range('E39').Value = '  +0.00%  '
Set-TextValue 'D40' '100.08'
$ws.Range('E40').Value = '  +11.14%  '
Set-TextValue 'D41' '0.0220'
$ws.Range('E41').Value = '  +5.88%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D42' '1.13'
$ws.Range('E42').Value = '  +2.69%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D43' '17.03'
$ws.Range('E43').Value = '  +4.70%  '
$ws.Range('E44').Value = '  -0.09%  '
Set-TextValue 'D45' '2.47'
$ws.Range('E45').Value = '  +3.61%  '
$ws.Range('D46').Value = '1.342.23'
$ws.Range('E46').Value = '  +0.24%  '
Set-TextValue 'D47' '2.44'
$ws.Range('E47').Value = '  +1.37%  '
Set-TextValue 'D48' '6.73'
$ws.Range('E48').Value = '  +3.40%  '
$ws.Range('E49').Value = '  -0.95%  '
Set-TextValue 'D50' '45.35'
$ws.Range('E50').Value = '  -7.63%  '
$ws.Range('E51').Value = '  +7.38%  '
